$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts B:F left to A:E, so the old F
# (taxon-count) column ends up as the new column E and the old A (style-1)
# data column is dropped.
$ws.Columns.Item(1).Delete()

# Fix the header text: MODEL_CONDITION -> MODELCONDITION (now in D1).
$ws.Range("D1").Value = "MODELCONDITION"
